$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.743.34"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.852.26"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -2.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4317"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07399"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8853"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.857.91"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.766"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.489"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07115"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009046"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.013"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "27.755.16"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.281"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "2.096.67"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.034"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.438"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08971"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.240"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.585"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.920"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.151"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.013"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05346"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01970"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.869"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1688"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.721"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4753"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.013"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.917"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
